# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-03-31 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-01 Monday", 2)

# Update the practice-problem table cells. Cells are addressed by
# (row, column) rather than by text Find/Replace because a couple of the
# old/new problem strings repeat elsewhere in the table, which would make a
# global text replace ambiguous.
$t = $d.Tables.Item(1)

$rows = @{
    1  = @("42÷8=5, 2", "70÷9=7, 7", "42÷6=7, 0", "62÷3=20, 2", "65÷2=32, 1")
    5  = @("44÷5=8, 4", "58÷2=29, 0", "89÷7=12, 5", "76÷2=38, 0", "18÷9=2, 0")
    9  = @("82÷6=13, 4", "86÷3=28, 2", "14÷2=7, 0", "89÷7=12, 5", "81÷9=9, 0")
    13 = @("66÷3=22, 0", "60÷3=20, 0", "44÷7=6, 2", "46÷3=15, 1", "83÷7=11, 6")
    17 = @("61÷3=20, 1", "58÷9=6, 4", "90÷3=30, 0", "22÷4=5, 2", "28÷3=9, 1")
}

foreach ($rowIndex in $rows.Keys) {
    $values = $rows[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Rows.Item($rowIndex).Cells.Item($col)
        $cell.Range.Text = $values[$col - 1]
    }
}
